$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: set resultado / profit
$ws.Range("G28").Value = "Acierto"
$ws.Range("H28").Value = 2.4

# Row 29: set resultado / profit
$ws.Range("G29").Value = "Fallo"
$ws.Range("H29").Value = -1

# Row 48: set resultado / profit
$ws.Range("G48").Value = "Fallo"
$ws.Range("H48").Value = -1

# Rows 43-51: event_id (column A) was stored as text; convert to a true number
$ws.Range("A43").Value = 14580793
$ws.Range("A44").Value = 14580346
$ws.Range("A45").Value = 14579620
$ws.Range("A46").Value = 14579621
$ws.Range("A47").Value = 14581354
$ws.Range("A48").Value = 14581358
$ws.Range("A49").Value = 14583788
$ws.Range("A50").Value = 14583789
$ws.Range("A51").Value = 14583779
